$d = $word.ActiveDocument
$rng = $d.Content
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body><w:p w14:paraId="1A9EE19E" w14:textId="3FBD1D28" w:rsidR="008C4F25" w:rsidRDefault="00D87EBA" w:rsidP="008C4F25"><w:r><w:t>Team Portfolio 465</w:t></w:r></w:p><w:p><w:r><w:t>About Us</w:t></w:r></w:p><w:p><w:r><w:t>Hello, and welcome to our group portfolio for Communications 465.  We have all been working together the entire semester on various projects we have been assigned for class</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve">  We have done tech breakdowns for local businesses, as well as studies on enterprise applications.  We have learned a lot about various ways to communicate with tech through COMM 465!</w:t></w:r></w:p><w:p w14:paraId="2B64790B" w14:textId="77777777" w:rsidR="0053079C" w:rsidRDefault="0053079C" w:rsidP="008C4F25"><w:r><w:t>Breakout Session #1</w:t></w:r></w:p><w:p w14:paraId="3C0C5381" w14:textId="77777777" w:rsidR="005232D3" w:rsidRDefault="0053079C" w:rsidP="008C4F25"><w:r><w:t>We did our work on the breakout session with bagel bakery</w:t></w:r><w:r w:rsidR="0007538B"><w:t xml:space="preserve"> as Susannah is their social media manager.  We wanted to see how we cou</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve">ld use technology to improve their business model.  We spoke with their manager and brainstormed to think of some ideas we could implement into their business.  We came up with a few changes to make to the website and implementing a new feature that the manager had been looking to </w:t></w:r><w:r w:rsidR="0061428E"><w:t>implement</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve"> for a while, an accurate bagel timer for when the bagels would be done.  To </w:t></w:r><w:r w:rsidR="0061428E"><w:t>implement</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve"> this Jeremy wrote a bit of </w:t></w:r><w:r w:rsidR="0061428E"><w:t>JavaScript</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve"> code himself to pull this off to the specifics that </w:t></w:r><w:r w:rsidR="0061428E"><w:t>the</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0061428E"><w:t>manager</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve"> wanted.  He also wrote some custom html and </w:t></w:r><w:r w:rsidR="003B57F1"><w:t>CSS</w:t></w:r><w:r w:rsidR="007D5076"><w:t xml:space="preserve"> code to properly embed the social medias into the site as well as just </w:t></w:r><w:r w:rsidR="0061428E"><w:t>spicing up the site design overall.  Here is a link to the site with the timer feature included.</w:t></w:r></w:p><w:p w14:paraId="65B2E4DA" w14:textId="77777777" w:rsidR="005232D3" w:rsidRDefault="005232D3" w:rsidP="008C4F25"><w:r><w:t xml:space="preserve">Breakout Session # 2 </w:t></w:r></w:p><w:p w14:paraId="1EF7EACD" w14:textId="7818B325" w:rsidR="007F64C7" w:rsidRDefault="005232D3" w:rsidP="008C4F25"><w:r><w:t xml:space="preserve">Here we analyzed an employee feedback tool called </w:t></w:r><w:r w:rsidR="001917CC"><w:t>Reflektive and</w:t></w:r><w:r><w:t xml:space="preserve"> gauged about the pros and cons of the software.  We found that for a </w:t></w:r><w:r w:rsidR="001917CC"><w:t>business-like</w:t></w:r><w:r><w:t xml:space="preserve"> bagel bakery</w:t></w:r><w:r w:rsidR="005E5190"><w:t xml:space="preserve"> that is a mom and Pop that this software would not really be a good fit. In </w:t></w:r><w:r w:rsidR="001917CC"><w:t>general,</w:t></w:r><w:r w:rsidR="005E5190"><w:t xml:space="preserve"> we also found that we had more grievances with the software than we had pros.  </w:t></w:r><w:r w:rsidR="005F07A3"><w:t xml:space="preserve">For the second part of the </w:t></w:r><w:r w:rsidR="001917CC"><w:t>project,</w:t></w:r><w:r w:rsidR="005F07A3"><w:t xml:space="preserve"> we all picked a job description for our careers of choice </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="005F07A3"><w:t>and also</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="005F07A3"><w:t xml:space="preserve"> talked about the job description in reference to our resumes.  This was a good exercise as many of us are about to be entering the </w:t></w:r><w:r w:rsidR="001917CC"><w:t>job-hunting</w:t></w:r><w:r w:rsidR="005F07A3"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="005F07A3"><w:lastRenderedPageBreak/><w:t>phase.</w:t></w:r><w:r w:rsidR="007F64C7"><w:br/><w:t>Wearables</w:t></w:r></w:p><w:p w14:paraId="2B1232EA" w14:textId="04B521CB" w:rsidR="007F64C7" w:rsidRDefault="007F64C7" w:rsidP="008C4F25"><w:r><w:t xml:space="preserve">Jeremy and Matt chose to do the topic of </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>wearable</w:t></w:r><w:r><w:t xml:space="preserve"> smart devices for their tech dive project.  </w:t></w:r><w:r w:rsidR="00C4179B"><w:t>They compared the current state of the market right now.  The three major players are Fitbit Wear OS (Android) and Apple</w:t></w:r><w:r w:rsidR="00FA1C66"><w:t xml:space="preserve">.  </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>Overall,</w:t></w:r><w:r w:rsidR="004F025A"><w:t xml:space="preserve"> they concluded that each category works best for different types of people.</w:t></w:r><w:r w:rsidR="000E1EEF"><w:t xml:space="preserve">  Fitbit works better for users that are </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>only</w:t></w:r><w:r w:rsidR="000E1EEF"><w:t xml:space="preserve"> focused on </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>exercise,</w:t></w:r><w:r w:rsidR="000E1EEF"><w:t xml:space="preserve"> don’t have a strong opinion about the apple or android ecosystem and / or aren’t very technologically literate.</w:t></w:r><w:r w:rsidR="007011AD"><w:t xml:space="preserve">  For Wear OS / android </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>watches,</w:t></w:r><w:r w:rsidR="007011AD"><w:t xml:space="preserve"> we found that this could be better for more </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>technologically</w:t></w:r><w:r w:rsidR="007011AD"><w:t xml:space="preserve"> savvy users that enjoy the android and google play </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>ecosystem</w:t></w:r><w:r w:rsidR="007011AD"><w:t xml:space="preserve"> and additionally enjoy the extra options the android market gives them.</w:t></w:r><w:r w:rsidR="003C5B55"><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>Finally,</w:t></w:r><w:r w:rsidR="003C5B55"><w:t xml:space="preserve"> the apple watch supports the average user the most, although pricey like most apple products, it supports full </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>integration</w:t></w:r><w:r w:rsidR="003C5B55"><w:t xml:space="preserve"> into the apple ecosystem and allows you to check messages. </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>Music</w:t></w:r><w:r w:rsidR="003C5B55"><w:t xml:space="preserve"> playlists calls </w:t></w:r><w:r w:rsidR="004E60EC"><w:t>etc.</w:t></w:r><w:r w:rsidR="003C5B55"><w:t xml:space="preserve"> all while also recording fitness stats.</w:t></w:r><w:r w:rsidR="004E60EC"><w:t xml:space="preserve">  We found that apple has the largest share in the market right now and that things will probably stay that way for a long time.</w:t></w:r></w:p><w:p w14:paraId="33DE275A" w14:textId="5CE8F679" w:rsidR="000A2934" w:rsidRDefault="000A2934" w:rsidP="008C4F25"><w:r><w:t>Bios</w:t></w:r></w:p><w:p w14:paraId="70447897" w14:textId="03D26172" w:rsidR="000A2934" w:rsidRDefault="000A2934" w:rsidP="008C4F25"><w:r><w:t>Jeremy</w:t></w:r></w:p><w:p w14:paraId="38F10CC9" w14:textId="7CFF4CF7" w:rsidR="000A2934" w:rsidRDefault="000A2934" w:rsidP="008C4F25"><w:r><w:t>I am</w:t></w:r><w:r w:rsidRPr="000A2934"><w:t xml:space="preserve"> finishing up my last semester at Salisbury University, I am receiving my </w:t></w:r><w:r w:rsidR="00DB5B80" w:rsidRPr="000A2934"><w:t>bachelor’s</w:t></w:r><w:r w:rsidRPr="000A2934"><w:t xml:space="preserve"> degree in Communications Media Production with a minor in Computer Science. I have quite the passion for everything Computer Science, networking, IT, programming, and everything in between! The last few years I have been teaching myself all that I can in my free time and in University Classes. I got interested in computer programming after I got an itch for game development after being inspired by Nintendo’s Earthbound game series. After this, I knew I wanted to learn more about </w:t></w:r><w:r w:rsidR="0053079C" w:rsidRPr="000A2934"><w:t>computer</w:t></w:r><w:r w:rsidRPr="000A2934"><w:t xml:space="preserve"> programming and tech, and this is what started my journey </w:t></w:r><w:r w:rsidRPr="000A2934"><w:lastRenderedPageBreak/><w:t xml:space="preserve">into pursuing CS and Software development. Immediately following my graduation in May 2022, I will be starting my </w:t></w:r><w:r w:rsidR="0053079C" w:rsidRPr="000A2934"><w:t>master’s</w:t></w:r><w:r w:rsidRPr="000A2934"><w:t xml:space="preserve"> program in Cybersecurity and Information Assurance at Western Governor’s University. My goal is to land a position as a fullstack Software Developer, after finishing graduate school.</w:t></w:r></w:p></w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)
Write-Output "Paragraphs: $($d.Paragraphs.Count)"
